$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, pushing existing rows 111-217 down to 112-218.
$ws.Rows("111").Insert()

# Populate the newly inserted row 111 with the new record's data.
$ws.Range("A111").Value = 7
$ws.Range("B111").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C111").Value = "Ñuble"
$ws.Range("D111").Value = 44658
$ws.Range("E111").Value = 16
$ws.Range("F111").Value = 100112003
$ws.Range("G111").Value = "Ajo"
$ws.Range("H111").Value = "Chino"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 60
$ws.Range("K111").Value = 19000
$ws.Range("L111").Value = 20000
$ws.Range("M111").Value = 19500
$ws.Range("N111").Value = "$/caja 10 kilos"
$ws.Range("O111").Value = "China"
$ws.Range("P111").Value = 1950
$ws.Range("Q111").Value = 10
$ws.Range("R111").Value = "Hortaliza"
